$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => D value (Price), E value (Volume 1h)
$updates = @(
    @{ Row = 2;  D = "63.862.37"; E = "  -0.20%  " },
    @{ Row = 3;  D = "3.146.07";  E = "  -0.08%  " },
    @{ Row = 4;  D = $null;       E = "  +0.09%  " },
    @{ Row = 5;  D = "586.18";    E = "  -0.81%  " },
    @{ Row = 6;  D = "146.05";    E = "  -1.04%  " },
    @{ Row = 8;  D = "3.133.91";  E = "  -0.24%  " },
    @{ Row = 9;  D = "0.527";     E = "  -1.73%  " },
    @{ Row = 10; D = "0.159";     E = "  -0.32%  " },
    @{ Row = 11; D = "5.78";      E = "  +0.80%  " },
    @{ Row = 12; D = "0.456";     E = "  -3.08%  " },
    @{ Row = 13; D = $null;       E = "  -2.91%  " },
    @{ Row = 14; D = "36.83";     E = "  +2.24%  " },
    @{ Row = 15; D = "3.665.14";  E = "  -0.18%  " },
    @{ Row = 16; D = $null;       E = "  -1.69%  " },
    @{ Row = 17; D = "63.622.18"; E = "  -0.47%  " },
    @{ Row = 18; D = "3.139.82";  E = "  -0.16%  " },
    @{ Row = 19; D = $null;       E = "  -2.01%  " },
    @{ Row = 20; D = "462.99";    E = "  -1.72%  " },
    @{ Row = 21; D = "14.29";     E = "  +0.35%  " },
    @{ Row = 22; D = $null;       E = "  -0.18%  " },
    @{ Row = 23; D = "7.40";      E = "  -2.24%  " },
    @{ Row = 24; D = "12.89";     E = "  -3.79%  " },
    @{ Row = 25; D = "80.90";     E = "  -2.07%  " },
    @{ Row = 26; D = $null;       E = "  +1.07%  " },
    @{ Row = 27; D = $null;       E = "  -0.16%  " },
    @{ Row = 28; D = "9.28";      E = "  +6.16%  " },
    @{ Row = 29; D = $null;       E = "  -1.39%  " },
    @{ Row = 30; D = $null;       E = "  -0.02%  " },
    @{ Row = 31; D = $null;       E = "  -0.76%  " },
    @{ Row = 32; D = "7.09";      E = "  +3.71%  " },
    @{ Row = 33; D = "26.84";     E = "  -1.06%  " },
    @{ Row = 34; D = "0.108";     E = "  -0.60%  " },
    @{ Row = 35; D = "0.0₃0844";  E = "  -3.52%  " },
    @{ Row = 36; D = $null;       E = "  -1.22%  " },
    @{ Row = 37; D = $null;       E = "  -4.58%  " },
    @{ Row = 38; D = "3.32";      E = "  -2.14%  " },
    @{ Row = 39; D = "5.99";      E = "  -2.80%  " },
    @{ Row = 40; D = "51.27";     E = "  +0.63%  " },
    @{ Row = 41; D = "436.58";    E = "  -2.98%  " },
    @{ Row = 42; D = $null;       E = "  +1.48%  " },
    @{ Row = 43; D = "2.910.83";  E = "  -0.54%  " },
    @{ Row = 44; D = $null;       E = "  -2.00%  " },
    @{ Row = 45; D = $null;       E = "  -1.38%  " },
    @{ Row = 46; D = $null;       E = "  -4.33%  " },
    @{ Row = 47; D = "37.47";     E = "  +8.41%  " },
    @{ Row = 48; D = "127.09";    E = "  +1.94%  " },
    @{ Row = 49; D = $null;       E = "  +0.00%  " },
    @{ Row = 50; D = $null;       E = "  -2.07%  " },
    @{ Row = 51; D = "24.00";     E = "  -3.82%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $cellE = $ws.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
}
